$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the per-category subtotal rows so each carries its section name ---

# Civilian
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Carryover nominations"
$ws.Range("A9").Value  = "     Civilian, Confirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

# Other Civilian
$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Carryover nominations"
$ws.Range("A15").Value = "     Other Civilian, Confirmed "
$ws.Range("A16").Value = "     Other Civilian, Withdrawn "
$ws.Range("A17").Value = "     Other Civilian, Returned to White House "

# Air Force
$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("A20").Value = "     Air Force, Carryover nominations"
$ws.Range("A21").Value = "     Air Force, Confirmed "
$ws.Range("A22").Value = "     Air Force, Withdrawn "
$ws.Range("A23").Value = "     Air Force, Returned to White House "

# Army
$ws.Range("A25").Value = "     Army, New nominations"
$ws.Range("A26").Value = "     Army, Carryover nominations"
$ws.Range("A27").Value = "     Army, Confirmed "
$ws.Range("A28").Value = "     Army, Withdrawn "
$ws.Range("A29").Value = "     Army, Returned to White House "

# Navy
$ws.Range("A31").Value = "     Navy, New nominations"
$ws.Range("A32").Value = "     Navy, Carryover nominations"
$ws.Range("A33").Value = "     Navy, Confirmed "
$ws.Range("A34").Value = "     Navy, Returned to White House "

# Marine Corps
$ws.Range("A36").Value = "     Marine Corps, New nominations"
$ws.Range("A37").Value = "     Marine Corps, Confirmed "

# --- Rework the "Summary" block (rows 38-44) into a 6-row totals block (38-43) ---
# Row 38 used to just be the "Summary" header; it now also carries the
# "Total new nominations" figure that used to live on row 40.
$ws.Range("A38").Value = "Total new nominations"
$ws.Range("B38").Value = 23803
$ws.Range("B38").NumberFormat = "#,##0"

$ws.Range("A39").Value = "Total carryover nominations"
$ws.Range("B39").Value = 667

$ws.Range("A40").Value = "Total confirmed "
$ws.Range("B40").Value = 24296

$ws.Range("A41").Value = "Total unconfirmed "
# This cell used to hold the "#,##0"-formatted 24296 figure; the new
# "Total unconfirmed" figure reuses the plain General-number formatting
# from its sibling totals (e.g. B39), so copy that formatting over first.
$ws.Range("B39").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B41").Value = 0

$ws.Range("A42").Value = "Total withdrawn "
$ws.Range("B42").Value = 21

$ws.Range("A43").Value = "Total returned to the White House "
$ws.Range("B43").Value = 153

# The old row 44 ("Total Returned to the White House ") is now redundant --
# its value was folded into row 43 above -- so remove it and let everything
# below it shift up.
$ws.Rows(44).Delete()
